$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 ("RM 232") entirely - shifts subsequent rows up
$ws.Rows.Item(26).Delete()

# Delete original row 28 ("SC 92"), which after the previous delete is now row 27
$ws.Rows.Item(27).Delete()

# Apply individual cell value corrections (new row numbers, after the deletions above)
$ws.Range("E5").ClearContents()
$ws.Range("E11").Value = -7.9
$ws.Range("C19").Value = 13.2
$ws.Range("E19").ClearContents()
$ws.Range("C21").ClearContents()
$ws.Range("C23").Value = 12.2
$ws.Range("E25").Value = -7.1
$ws.Range("C27").ClearContents()
$ws.Range("E29").ClearContents()
$ws.Range("C33").Value = 10.4

Write-Output "done"
